$wb = $excel.ActiveWorkbook

# Update cached market-price / profit values per scheduled runner refresh.
# Each block updates one row of a sheet; values come from the latest price pull.

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 50591.418
$ws.Range("I32").Value = 4499.75
$ws.Range("J32").Value = 73637.25
$ws.Range("K32").Value = 4499.75
$ws.Range("L32").Value = 73637.25
$ws.Range("M32").Value = -4173.75
$ws.Range("N32").Value = -74289.25

# Row 74
$ws.Range("H74").Value = 5749.75
$ws.Range("I74").Value = 5999.6665
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 5999.6665
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -5063.6665
$ws.Range("N74").Value = -6872

# Row 76
$ws.Range("H76").Value = 5814.2856
$ws.Range("I76").Value = 3900
$ws.Range("K76").Value = 3900
$ws.Range("M76").Value = -3585

# Row 77
$ws.Range("H77").Value = 5749.75
$ws.Range("I77").Value = 5999.6665
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 29998.3325
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -25318.3325
$ws.Range("N77").Value = -34360

# Row 79
$ws.Range("H79").Value = 5814.2856
$ws.Range("I79").Value = 3900
$ws.Range("K79").Value = 3900
$ws.Range("M79").Value = -2808

# Row 103
$ws.Range("H103").Value = 689.4
$ws.Range("I103").Value = 749.5
$ws.Range("J103").Value = 449
$ws.Range("K103").Value = 2248.5
$ws.Range("L103").Value = 1347
$ws.Range("M103").Value = -1662.5
$ws.Range("N103").Value = -2519

# Row 132
$ws.Range("H132").Value = 25642238
$ws.Range("I132").Value = 28572442
$ws.Range("J132").Value = 2944.5
$ws.Range("K132").Value = 85717326
$ws.Range("L132").Value = 8833.5
$ws.Range("M132").Value = -85714796
$ws.Range("N132").Value = -13893.5

# Row 135
$ws.Range("H135").Value = 2299.6667
$ws.Range("I135").Value = 2112.0732
$ws.Range("J135").Value = 9991
$ws.Range("K135").Value = 19008.6588
$ws.Range("L135").Value = 89919
$ws.Range("M135").Value = -16473.6588
$ws.Range("N135").Value = -94989

# Row 137
$ws.Range("H137").Value = 3150.6667
$ws.Range("I137").Value = 2914.1333
$ws.Range("J137").Value = 4333.3335
$ws.Range("K137").Value = 8742.3999
$ws.Range("L137").Value = 13000.0005
$ws.Range("M137").Value = -6192.3999
$ws.Range("N137").Value = -18100.0005

# Row 138
$ws.Range("H138").Value = 3007.052
$ws.Range("I138").Value = 1720.9412
$ws.Range("J138").Value = 4023.9768
$ws.Range("K138").Value = 5162.8236
$ws.Range("L138").Value = 12071.9304
$ws.Range("M138").Value = -22.82359999999971
$ws.Range("N138").Value = -22351.9304

# Row 141
$ws.Range("H141").Value = 2180.4285
$ws.Range("I141").Value = 2180.4285
$ws.Range("K141").Value = 6541.2855
$ws.Range("M141").Value = -1361.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6546.8125
$ws.Range("I2").Value = 892.8461
$ws.Range("K2").Value = 892.8461
$ws.Range("M2").Value = -779.8461

# Row 32
$ws.Range("H32").Value = 2634.625
$ws.Range("I32").Value = 2284.1614
$ws.Range("J32").Value = 13499
$ws.Range("K32").Value = 2284.1614
$ws.Range("L32").Value = 13499
$ws.Range("M32").Value = -1997.1614
$ws.Range("N32").Value = -14073

# Row 61
$ws.Range("H61").Value = 1899.6666
$ws.Range("I61").Value = 1594.65
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 1594.65
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -1382.65
$ws.Range("N61").Value = -8424

# Row 74
$ws.Range("H74").Value = 3620.5625
$ws.Range("I74").Value = 2367.75
$ws.Range("J74").Value = 4873.375
$ws.Range("K74").Value = 2367.75
$ws.Range("L74").Value = 4873.375
$ws.Range("M74").Value = -1493.75
$ws.Range("N74").Value = -6621.375

# Row 77
$ws.Range("H77").Value = 3620.5625
$ws.Range("I77").Value = 2367.75
$ws.Range("J77").Value = 4873.375
$ws.Range("K77").Value = 11838.75
$ws.Range("L77").Value = 24366.875
$ws.Range("M77").Value = -7470.75
$ws.Range("N77").Value = -33102.875

# Row 116
$ws.Range("H116").Value = 6546.8125
$ws.Range("I116").Value = 892.8461
$ws.Range("K116").Value = 892.8461
$ws.Range("M116").Value = 1401.1539

# Row 132
$ws.Range("H132").Value = 8396.962
$ws.Range("I132").Value = 5186.3267
$ws.Range("J132").Value = 60837.332
$ws.Range("K132").Value = 15558.9801
$ws.Range("L132").Value = 182511.996
$ws.Range("M132").Value = -13028.9801
$ws.Range("N132").Value = -187571.996

# Row 136
$ws.Range("H136").Value = 1899.6666
$ws.Range("I136").Value = 1594.65
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 4783.950000000001
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -2233.950000000001
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6546.8125
$ws.Range("I3").Value = 892.8461
$ws.Range("K3").Value = 892.8461
$ws.Range("M3").Value = -778.8461

# Row 20
$ws.Range("H20").Value = 52633428
$ws.Range("I20").Value = 90911140
$ws.Range("J20").Value = 1573.875
$ws.Range("K20").Value = 90911140
$ws.Range("L20").Value = 1573.875
$ws.Range("M20").Value = -90910893
$ws.Range("N20").Value = -2067.875

# Row 134
$ws.Range("H134").Value = 3083.5103
$ws.Range("I134").Value = 1916.4634
$ws.Range("K134").Value = 5749.3902
$ws.Range("M134").Value = -3214.3902

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3449.5
$ws.Range("I31").Value = 3167.0588
$ws.Range("J31").Value = 5050
$ws.Range("K31").Value = 3167.0588
$ws.Range("L31").Value = 5050
$ws.Range("M31").Value = -2872.0588
$ws.Range("N31").Value = -5640

# Row 34
$ws.Range("H34").Value = 3449.5
$ws.Range("I34").Value = 3167.0588
$ws.Range("J34").Value = 5050
$ws.Range("K34").Value = 3167.0588
$ws.Range("L34").Value = 5050
$ws.Range("M34").Value = -2965.0588
$ws.Range("N34").Value = -5454

# Row 58
$ws.Range("H58").Value = 2912.2068
$ws.Range("I58").Value = 2976.9285
$ws.Range("K58").Value = 2976.9285
$ws.Range("M58").Value = -2773.9285

# Row 107
$ws.Range("H107").Value = 807.3913
$ws.Range("I107").Value = 440.57144
$ws.Range("J107").Value = 967.875
$ws.Range("K107").Value = 440.57144
$ws.Range("L107").Value = 967.875
$ws.Range("M107").Value = 1479.42856
$ws.Range("N107").Value = -4807.875

# Row 132
$ws.Range("H132").Value = 1319.4839
$ws.Range("I132").Value = 1246.5834
$ws.Range("J132").Value = 1569.4286
$ws.Range("K132").Value = 3739.7502
$ws.Range("L132").Value = 4708.2858
$ws.Range("M132").Value = -1209.7502
$ws.Range("N132").Value = -9768.2858

# Row 134
$ws.Range("H134").Value = 2068.9565
$ws.Range("I134").Value = 2068.9565
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6206.869499999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3671.869499999999
$ws.Range("N134").ClearContents()

# Row 136
$ws.Range("H136").Value = 2912.2068
$ws.Range("I136").Value = 2976.9285
$ws.Range("K136").Value = 8930.7855
$ws.Range("M136").Value = -6380.7855

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 2181.6667
$ws.Range("I8").Value = 2181.6667
$ws.Range("K8").Value = 6545.000100000001
$ws.Range("M8").Value = -6406.000100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1000.8077
$ws.Range("I132").Value = 984.2083
$ws.Range("K132").Value = 2952.6249
$ws.Range("M132").Value = -422.6248999999998

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 132
$ws.Range("H132").Value = 6145.0386
$ws.Range("I132").Value = 4418.6
$ws.Range("J132").Value = 8499.272000000001
$ws.Range("K132").Value = 13255.8
$ws.Range("L132").Value = 25497.816
$ws.Range("M132").Value = -10725.8
$ws.Range("N132").Value = -30557.816

# Row 136
$ws.Range("H136").Value = 5002.077
$ws.Range("I136").Value = 4483.6665
$ws.Range("J136").Value = 7179.4
$ws.Range("K136").Value = 13450.9995
$ws.Range("L136").Value = 21538.2
$ws.Range("M136").Value = -10900.9995
$ws.Range("N136").Value = -26638.2

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1378
$ws.Range("I132").Value = 833.4211
$ws.Range("K132").Value = 2500.2633
$ws.Range("M132").Value = 29.73669999999993

# Row 136
$ws.Range("H136").Value = 5938.2646
$ws.Range("I136").Value = 7188.885
$ws.Range("J136").Value = 1873.75
$ws.Range("K136").Value = 21566.655
$ws.Range("L136").Value = 5621.25
$ws.Range("M136").Value = -19016.655
$ws.Range("N136").Value = -10721.25
